$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates. Column D (Price) values are numeric-looking strings
# that must remain stored as text (matching the source inlineStr cells),
# so we prefix with an apostrophe to force text entry, then reset the
# cell style back to Normal so no visible quote-prefix / style drift remains.

$c = $ws.Range("D2")
$c.Value = "'" + '29.392.08'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -0.47%  '
$c = $ws.Range("D3")
$c.Value = "'" + '1.847.03'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.08%  '
$c = $ws.Range("D5")
$c.Value = "'" + '240.20'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.02%  '
$c = $ws.Range("D6")
$c.Value = "'" + '0.6317'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.01%  '
$c = $ws.Range("D7")
$c.Value = "'" + '1.0000'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +0.01%  '
$c = $ws.Range("D9")
$c.Value = "'" + '0.2966'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.28%  '
$c = $ws.Range("D10")
$c.Value = "'" + '24.60'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +1.24%  '
$c = $ws.Range("D11")
$c.Value = "'" + '0.07729'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.68%  '
$c = $ws.Range("D12")
$c.Value = "'" + '1.845.06'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -1.57%  '
$c = $ws.Range("D13")
$c.Value = "'" + '5.000'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -0.45%  '
$c = $ws.Range("D14")
$c.Value = "'" + '0.6856'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +0.09%  '
$c = $ws.Range("D15")
$c.Value = "'" + '0.00001001'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +1.71%  '
$c = $ws.Range("D16")
$c.Value = "'" + '83.12'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.73%  '
$ws.Range("E17").Value = '  -0.40%  '
$c = $ws.Range("D18")
$c.Value = "'" + '29.418.01'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -0.54%  '
$c = $ws.Range("D19")
$c.Value = "'" + '229.90'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.81%  '
$ws.Range("E20").Value = '  -0.26%  '
$c = $ws.Range("D21")
$c.Value = "'" + '0.9993'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.10%  '
$c = $ws.Range("D22")
$c.Value = "'" + '7.578'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.37%  '
$ws.Range("E23").Value = '  -0.03%  '
$c = $ws.Range("D24")
$c.Value = "'" + '157.13'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.90%  '
$c = $ws.Range("D25")
$c.Value = "'" + '0.1400'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.98%  '
$c = $ws.Range("D26")
$c.Value = "'" + '8.384'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.52%  '
$ws.Range("E27").Value = '  -0.18%  '
$ws.Range("E28").Value = '  -1.02%  '
$c = $ws.Range("D29")
$c.Value = "'" + '0.05732'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -1.82%  '
$c = $ws.Range("D30")
$c.Value = "'" + '1.252'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -2.33%  '
$ws.Range("E31").Value = '  +0.48%  '
$ws.Range("E32").Value = '  -0.16%  '
$ws.Range("E33").Value = '  -2.33%  '
$c = $ws.Range("D35")
$c.Value = "'" + '0.7178'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +0.28%  '
$c = $ws.Range("D36")
$c.Value = "'" + '2.593'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +0.07%  '
$c = $ws.Range("D37")
$c.Value = "'" + '1.250.64'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +0.85%  '
$c = $ws.Range("D38")
$c.Value = "'" + '0.01819'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +2.56%  '
$c = $ws.Range("D39")
$c.Value = "'" + '2.784'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.56%  '
$c = $ws.Range("D40")
$c.Value = "'" + '6.214'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +1.19%  '
$c = $ws.Range("D41")
$c.Value = "'" + '0.9074'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -0.74%  '
$c = $ws.Range("D42")
$c.Value = "'" + '1.001'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("B43").Value = 'RocketPoolETH'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$c = $ws.Range("D43")
$c.Value = "'" + '1.999.31'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -2.07%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Range("D44")
$c.Value = "'" + '101.75'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -0.16%  '
$c = $ws.Range("D45")
$c.Value = "'" + '66.48'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -1.52%  '
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range("D46")
$c.Value = "'" + '7.063'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -2.78%  '
$c = $ws.Range("D47")
$c.Value = "'" + '9.166'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +0.19%  '
$ws.Range("B48").Value = 'TheSandbox'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Range("D48")
$c.Value = "'" + '0.4031'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.06%  '
$c = $ws.Range("D49")
$c.Value = "'" + '1.711'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +0.62%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range("D50")
$c.Value = "'" + '0.1131'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +0.99%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range("D51")
$c.Value = "'" + '0.05743'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.13%  '
